# Generate Report for Handoff
#
# This applies a "handoff" status refresh to the localization-status workbook:
#   - The row describing the now-obsolete "ea39ce26-...md" file is removed
#     from every sheet (Overview, zh-cn, de-de), along with its hyperlinks.
#   - The remaining row for "91bbbb67-...md" is updated from
#     "Handed back: in sync with en-US" to "Ready for handoff", with fresh
#     handoff timestamps.

function Remove-RowHyperlinks {
    param($ws, [int]$row)

    # Deleting a Hyperlink object can invalidate other cached Hyperlink
    # references in the same collection, so re-scan the live collection
    # after every deletion instead of deleting from a pre-built list.
    $changed = $true
    while ($changed) {
        $changed = $false
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Row -eq $row) {
                $h.Delete()
                $changed = $true
                break
            }
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Remove-RowHyperlinks $ws1 3
$ws1.Rows("3:3").Delete()

$ws1.Range("B2").Value2 = "Ready for handoff"
$ws1.Range("C2").Value2 = "Ready for handoff"
$ws1.Range("D2").Value2 = "2016-03-24 10:23:19"

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

Remove-RowHyperlinks $ws2 3
$ws2.Rows("3:3").Delete()

$ws2.Range("C2").Value2 = "Ready for handoff"
$ws2.Range("E2").Value2 = "2016-03-24 10:23:10"

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

Remove-RowHyperlinks $ws3 3
$ws3.Rows("3:3").Delete()

$ws3.Range("C2").Value2 = "Ready for handoff"
$ws3.Range("E2").Value2 = "2016-03-24 10:23:19"
